$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLoginTest")

$ws.Range("A4").Value = "peter"
$ws.Range("B4").Value = "peter123"
$ws.Range("C4").Value = "Invalid credential"

$ws.Range("E5").Select()
